# Update FlashScore odds data for 2024-12-10 per upstream refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 3
$ws.Range("G3").Value = 1.17
$ws.Range("H3").Value = 7
$ws.Range("J3").Value = 1.53
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 10
$ws.Range("O3").Value = 1.11
$ws.Range("P3").Value = 6.5
$ws.Range("Q3").Value = 1.4
$ws.Range("R3").Value = 2.88
$ws.Range("U3").Value = 2
$ws.Range("V3").Value = 1.73
$ws.Range("Y3").Value = 10
$ws.Range("AB3").Value = 29
$ws.Range("AC3").Value = 21
$ws.Range("AD3").Value = 15
$ws.Range("AE3").Value = 26
$ws.Range("AF3").Value = 67
$ws.Range("AG3").Value = 351
$ws.Range("AH3").Value = 34
$ws.Range("AI3").Value = 67
$ws.Range("AJ3").Value = 34
$ws.Range("AL3").Value = 81
$ws.Range("AM3").Value = 67
$ws.Range("AN3").Value = 3.25
$ws.Range("AO3").Value = 5
$ws.Range("AQ3").Value = 11
$ws.Range("AX3").Value = 13
$ws.Range("AY3").Value = 51
$ws.Range("BA3").Value = 251
$ws.Range("BB3").Value = 251
$ws.Range("BC3").Value = 351

# Row 4
$ws.Range("G4").Value = 2
$ws.Range("I4").Value = 3.5
$ws.Range("J4").Value = 2.63
$ws.Range("L4").Value = 4
$ws.Range("O4").Value = 1.25
$ws.Range("P4").Value = 3.75
$ws.Range("W4").Value = 8.5
$ws.Range("X4").Value = 10
$ws.Range("Z4").Value = 17
$ws.Range("AA4").Value = 15
$ws.Range("AH4").Value = 12
$ws.Range("AI4").Value = 19
$ws.Range("AK4").Value = 41
$ws.Range("AL4").Value = 26
$ws.Range("AM4").Value = 34
$ws.Range("AN4").Value = 4
$ws.Range("AP4").Value = 19
$ws.Range("AQ4").Value = 34
$ws.Range("AY4").Value = 19
$ws.Range("AZ4").Value = 26
$ws.Range("BB4").Value = 81

# Row 6
$ws.Range("AW6").Value = 126
$ws.Range("BD6").Value = 126

# Row 7
$ws.Range("O7").Value = 1.22
$ws.Range("P7").Value = 4.33
$ws.Range("AR7").Value = 81

# Row 8
$ws.Range("G8").Value = 3.2
$ws.Range("I8").Value = 2.2
$ws.Range("W8").Value = 10
$ws.Range("Y8").Value = 11
$ws.Range("AA8").Value = 23
$ws.Range("AH8").Value = 8.5
$ws.Range("AN8").Value = 5
$ws.Range("AO8").Value = 17

# Row 9
$ws.Range("Q9").Value = 2.2
$ws.Range("R9").Value = 1.67

# Row 11
$ws.Range("G11").Value = 1.67
$ws.Range("H11").Value = 3.8
$ws.Range("J11").Value = 2.25
$ws.Range("K11").Value = 2.3
$ws.Range("M11").Value = 1.04
$ws.Range("N11").Value = 13
$ws.Range("O11").Value = 1.25
$ws.Range("P11").Value = 4
$ws.Range("Q11").Value = 1.8
$ws.Range("R11").Value = 2
$ws.Range("S11").Value = 1.36
$ws.Range("T11").Value = 3
$ws.Range("U11").Value = 1.73
$ws.Range("V11").Value = 2
$ws.Range("W11").Value = 7.5
$ws.Range("X11").Value = 8.5
$ws.Range("AB11").Value = 23
$ws.Range("AC11").Value = 12
$ws.Range("AD11").Value = 7.5
$ws.Range("AE11").Value = 15
$ws.Range("AG11").Value = 201
$ws.Range("AH11").Value = 15
$ws.Range("AI11").Value = 26
$ws.Range("AO11").Value = 8.5
$ws.Range("AP11").Value = 19
$ws.Range("AS11").Value = 126
$ws.Range("AT11").Value = 3
$ws.Range("AU11").Value = 8
$ws.Range("AZ11").Value = 29
$ws.Range("BC11").Value = 201

# Row 14
$ws.Range("G14").Value = 2.4
$ws.Range("H14").Value = 3.15
$ws.Range("I14").Value = 2.72
$ws.Range("J14").Value = 3
$ws.Range("K14").Value = 2.12
$ws.Range("L14").Value = 3.25
$ws.Range("N14").Value = 7.9
$ws.Range("P14").Value = 3.75
$ws.Range("S14").Value = 1.36
$ws.Range("T14").Value = 2.92
$ws.Range("U14").Value = 1.55
$ws.Range("W14").Value = 9.75
$ws.Range("X14").Value = 13.5
$ws.Range("Y14").Value = 9
$ws.Range("Z14").Value = 27
$ws.Range("AA14").Value = 18.5
$ws.Range("AC14").Value = 7.9
$ws.Range("AD14").Value = 6.4
$ws.Range("AH14").Value = 11
$ws.Range("AI14").Value = 16.5
$ws.Range("AJ14").Value = 9.75
$ws.Range("AK14").Value = 35
$ws.Range("AL14").Value = 21
$ws.Range("AM14").Value = 24
$ws.Range("AN14").Value = 4.6
$ws.Range("AQ14").Value = 50
$ws.Range("AR14").Value = 75
$ws.Range("AS14").Value = 200
$ws.Range("AT14").Value = 2.92
$ws.Range("AX14").Value = 4.9
$ws.Range("AY14").Value = 14.5
$ws.Range("AZ14").Value = 19
$ws.Range("BB14").Value = 80
$ws.Range("BC14").Value = 200
